$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "run" column shifted from 0 to 1 for the first four data rows
$ws.Range("A2:A5").Value = "1"

# Give column B (alpha_name) an explicit width, as after a manual fit
$ws.Columns.Item(2).ColumnWidth = 14

# Leave the selection where the user last left it before saving
$ws.Range("J19").Select() | Out-Null
